$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: risk "HRD Emplloyees not cooperating due to unhapiness" ---
# Probability mark moves from F5 (L) to E5 (VL); new mark gets underlined emphasis.
$ws.Range("F5").ClearContents()
$ws.Range("E5").Value = "X"
$ws.Range("E5").Font.Underline = $true

# --- Row 6: risk "HRD Delegate leaving the project" ---
# Probability mark moves from G6 (M) to E6 (VL).
$ws.Range("G6").ClearContents()
$ws.Range("E6").Value = "X"

# --- Row 8: risk "Technical Coordenator leaving the project" ---
# Probability mark moves from G8 (M) to F8 (L); Consequence mark moves from N8 (VH) to L8 (M).
$ws.Range("G8").Clear()
$ws.Range("N8").Clear()
$ws.Range("F8").Value = "X"
$ws.Range("F8").Font.Underline = $true
$ws.Range("L8").Value = "X"

# --- View refresh to reflect the updated status ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("R10").Select() | Out-Null
